# PurchaseList.xlsx — all resistors & capacitors from 0603 to 0402,
# new eurocircuits basket B2206781

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "report created" time stamp (F24): 17:43 -> 11:49
# ---------------------------------------------------------------------------
$ws.Range("F24").Value = "11:49"

# ---------------------------------------------------------------------------
# Footprints (column D): 0603 packages -> 0402 packages
# ---------------------------------------------------------------------------
# Capacitors C1,C4,C5 / C2,C3 / C6,C7 / C8,C10
$ws.Range("D2:D5").Value = "0402_CAP"

# Resistors R1,R3,R4,R9,R14 / R2 / R5,R7 / R6,R8 / R10 / R12 / R15
$ws.Range("D12:D18").Value = "0402_res"
# R5,R7 and R15 use the upper-case footprint variant
$ws.Range("D14").Value = "0402_RES"
$ws.Range("D18").Value = "0402_RES"

# ---------------------------------------------------------------------------
# Row 2 - C1, C4, C5 - 100nF
# ---------------------------------------------------------------------------
$ws.Range("F2").Value = "Wurth Electronics"
$ws.Range("G2").Value = "885012105016"
$ws.Range("J2").Value = 40996
$ws.Range("K2").Value = "710-885012105016"
$ws.Range("P2").Value = 0.01885

# ---------------------------------------------------------------------------
# Row 3 - C2, C3 - 18pF
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = "CAPACITOR"
$ws.Range("F3").Value = "Walsin Technologies"
$ws.Range("G3").Value = "0402N180F500CT"
$ws.Range("J3").Value = 5900
$ws.Range("K3").Value = "791-0402N180F500CT"
$ws.Range("P3").Value = 0.03591

# ---------------------------------------------------------------------------
# Row 4 - C6, C7 - 33pF
# ---------------------------------------------------------------------------
$ws.Range("E4").Value = "CAPACITOR"
$ws.Range("F4").Value = "Wurth Electronics"
$ws.Range("G4").Value = "885012005058"
$ws.Range("J4").Value = 34052
$ws.Range("K4").Value = "710-885012005058"
$ws.Range("P4").Value = 0.01257

# ---------------------------------------------------------------------------
# Row 5 - C8, C10 - 1uF/MLCC
# ---------------------------------------------------------------------------
$ws.Range("E5").Value = "CAPACITOR"
$ws.Range("F5").Value = "Taiyo Yuden"
$ws.Range("G5").Value = "JMK105BJ105KP-F"
$ws.Range("J5").Value = 25442
$ws.Range("K5").Value = "963-JMK105BJ105KP-F"
$ws.Range("P5").Value = 0.02873

# ---------------------------------------------------------------------------
# Row 6 - D1 - LED_RED_0603 (stock / price refresh only)
# ---------------------------------------------------------------------------
$ws.Range("F6").Value = "Osram Opto"
$ws.Range("G6").Value = "LSQ976-NR-1"
$ws.Range("J6").Value = 68575
$ws.Range("K6").Value = "720-LSQ976-NR-1"
$ws.Range("P6").Value = 0.06822

# ---------------------------------------------------------------------------
# Row 7 - D2 - LED_GRN_0603 (stock / price refresh only)
# ---------------------------------------------------------------------------
$ws.Range("F7").Value = "Kingbright"
$ws.Range("G7").Value = "APT1608SGC"
$ws.Range("J7").Value = 43490
$ws.Range("K7").Value = "604-APT1608SGC"
$ws.Range("P7").Value = 0.05027

# ---------------------------------------------------------------------------
# Row 8 - D3 - LED_BLU_0603 (stock / price refresh only)
# ---------------------------------------------------------------------------
$ws.Range("F8").Value = "Osram Opto"
$ws.Range("G8").Value = "LB Q39G-L2OO-35-1"
$ws.Range("J8").Value = 58553
$ws.Range("K8").Value = "720-LBQ39GL2N2351"
$ws.Range("P8").Value = 0.10862

# ---------------------------------------------------------------------------
# Row 9 - D5 - NSR20F30NXT5G (stock / price refresh only)
# ---------------------------------------------------------------------------
$ws.Range("E9").Value = "SCHOTTKY DIODE 0603"
$ws.Range("F9").Value = "ON Semiconductor"
$ws.Range("K9").Value = "863-NSR20F30NXT5G"
$ws.Range("P9").Value = 0.18582

# ---------------------------------------------------------------------------
# Row 10 - JP2 - USB-MINI-B (stock / price refresh only)
# ---------------------------------------------------------------------------
$ws.Range("F10").Value = "Hirose"
$ws.Range("G10").Value = "UX60A-MB-5ST"
$ws.Range("J10").Value = 21568
$ws.Range("K10").Value = "798-UX60A-MB-5ST"
$ws.Range("P10").Value = 0.61312

# ---------------------------------------------------------------------------
# Row 11 - Q1 - BSS84,215 (stock / price refresh only)
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = "POWER MOSFET-P SOT23"
$ws.Range("F11").Value = "Diodes"
$ws.Range("G11").Value = "BSS84W-7-F"
$ws.Range("J11").Value = 97989
$ws.Range("K11").Value = "621-BSS84W-F"
$ws.Range("P11").Value = 0.10234

# ---------------------------------------------------------------------------
# Row 12 - R1, R3, R4, R9, R14 - 10K
# ---------------------------------------------------------------------------
$ws.Range("E12").Value = "RESISTOR"
$ws.Range("F12").Value = "Vishay"
$ws.Range("G12").Value = "CRCW040210K0FKEDC"
$ws.Range("J12").Value = 1868158
$ws.Range("K12").Value = "71-CRCW040210K0FKEDC"
$ws.Range("P12").Value = 0.00449

# ---------------------------------------------------------------------------
# Row 13 - R2 - 680R
# ---------------------------------------------------------------------------
$ws.Range("E13").Value = "RESISTOR"
$ws.Range("F13").Value = "Vishay"
$ws.Range("G13").Value = "CRCW0402680RFKEDC"
$ws.Range("J13").Value = 69395
$ws.Range("K13").Value = "71-CRCW0402680RFKEDC"
$ws.Range("P13").Value = 0.01167

# ---------------------------------------------------------------------------
# Row 14 - R5, R7 - 470R
# ---------------------------------------------------------------------------
$ws.Range("E14").Value = "RESISTOR"
$ws.Range("F14").Value = "Vishay Semiconductors"
$ws.Range("G14").Value = "CRCW0402470RFKEDC"
$ws.Range("J14").Value = 84527
$ws.Range("K14").Value = "71-CRCW0402470RFKEDC"
$ws.Range("P14").Value = 0.01167

# ---------------------------------------------------------------------------
# Row 15 - R6, R8 - 33R
# ---------------------------------------------------------------------------
$ws.Range("E15").Value = "RESISTOR"
$ws.Range("F15").Value = "Vishay Semiconductors"
$ws.Range("G15").Value = "CRCW040233R0FKEDC"
$ws.Range("J15").Value = 135537
$ws.Range("K15").Value = "71-CRCW040233R0FKEDC"
$ws.Range("P15").Value = 0.01167

# ---------------------------------------------------------------------------
# Row 16 - R10 - 1.5K
# ---------------------------------------------------------------------------
$ws.Range("E16").Value = "RESISTOR"
$ws.Range("F16").Value = "Vishay"
$ws.Range("G16").Value = "CRCW04021K50FKEDC"
$ws.Range("J16").Value = 99644
$ws.Range("K16").Value = "71-CRCW04021K50FKEDC"
$ws.Range("P16").Value = 0.01167

# ---------------------------------------------------------------------------
# Row 17 - R12 - 0R
# ---------------------------------------------------------------------------
$ws.Range("E17").Value = "RESISTOR"
$ws.Range("F17").Value = "Vishay"
$ws.Range("G17").Value = "CRCW04020000Z0EDC"
$ws.Range("J17").Value = 900334
$ws.Range("K17").Value = "71-CRCW04020000Z0EDC"
$ws.Range("P17").Value = 0.00539

# ---------------------------------------------------------------------------
# Row 18 - R15 - 18K
# ---------------------------------------------------------------------------
$ws.Range("E18").Value = "RESISTOR"
$ws.Range("F18").Value = "Vishay"
$ws.Range("G18").Value = "CRCW040218K0FKED"
$ws.Range("J18").Value = 18062
$ws.Range("K18").Value = "71-CRCW0402-18K-E3"
$ws.Range("P18").Value = 0.02424

# ---------------------------------------------------------------------------
# Row 19 - SW3, SW4 - EVQ-P2202M (stock / price refresh only)
# ---------------------------------------------------------------------------
$ws.Range("K19").Value = "667-EVQ-P2202M"
$ws.Range("P19").Value = 0.48744

# ---------------------------------------------------------------------------
# Row 20 - U1 - LPC11U35FHI33/501 (stock / price refresh only)
# ---------------------------------------------------------------------------
$ws.Range("K20").Value = "771-LPC11U35FHI33501"
$ws.Range("P20").Value = 2.98

# ---------------------------------------------------------------------------
# Row 21 - U2 - TPS78233 (stock / price refresh only)
# ---------------------------------------------------------------------------
$ws.Range("K21").Value = "595-TPS78233DDCR"
$ws.Range("P21").Value = 0.3878

# ---------------------------------------------------------------------------
# Row 22 - Y1 - 12MHz (stock / price refresh only)
# ---------------------------------------------------------------------------
$ws.Range("J22").Value = 930
$ws.Range("K22").Value = "717-8Z-12.000MAAJ-T"
$ws.Range("P22").Value = 1.01

# ---------------------------------------------------------------------------
# Row heights: three rows now wrap their longer Supplier Part Number text
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 21
$ws.Rows.Item(14).RowHeight = 21
$ws.Rows.Item(15).RowHeight = 21
